$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2563.3809
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2563.3809
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 7690.1427
$ws.Range("N17").Value = -8026.1427
$ws.Range("H33").Value = 272.8095
$ws.Range("I33").Value = 98.86667
$ws.Range("J33").Value = 707.6667
$ws.Range("K33").Value = 98.86667
$ws.Range("L33").Value = 707.6667
$ws.Range("M33").Value = 130.13333
$ws.Range("N33").Value = -1165.6667
$ws.Range("H40").Value = 1879.8
$ws.Range("I40").Value = 1879.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1879.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1704.8
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 51454.855
$ws.Range("I70").Value = 1735.5714
$ws.Range("J70").Value = 101174.14
$ws.Range("K70").Value = 5206.7142
$ws.Range("L70").Value = 303522.42
$ws.Range("M70").Value = -4936.7142
$ws.Range("N70").Value = -304062.42
$ws.Range("H73").Value = 51454.855
$ws.Range("I73").Value = 1735.5714
$ws.Range("J73").Value = 101174.14
$ws.Range("K73").Value = 5206.7142
$ws.Range("L73").Value = 303522.42
$ws.Range("M73").Value = -4270.7142
$ws.Range("N73").Value = -305394.42
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
$ws.Range("H132").Value = 1768.4694
$ws.Range("I132").Value = 1681
$ws.Range("J132").Value = 2395.3333
$ws.Range("K132").Value = 5043
$ws.Range("L132").Value = 7185.999899999999
$ws.Range("M132").Value = -2513
$ws.Range("N132").Value = -12245.9999
$ws.Range("H138").Value = 5709.8
$ws.Range("I138").Value = 7568.8887
$ws.Range("J138").Value = 5170.0645
$ws.Range("K138").Value = 22706.6661
$ws.Range("L138").Value = 15510.1935
$ws.Range("M138").Value = -17566.6661
$ws.Range("N138").Value = -25790.1935

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2126.3333
$ws.Range("I2").Value = 2242.125
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 2242.125
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -2129.125
$ws.Range("N2").Value = -1426
$ws.Range("H32").Value = 17954.508
$ws.Range("I32").Value = 7976.5527
$ws.Range("J32").Value = 30185.549
$ws.Range("K32").Value = 7976.5527
$ws.Range("L32").Value = 30185.549
$ws.Range("M32").Value = -7689.5527
$ws.Range("N32").Value = -30759.549
$ws.Range("H97").Value = 596.6667
$ws.Range("I97").Value = 583.3333
$ws.Range("J97").Value = 676.6667
$ws.Range("K97").Value = 583.3333
$ws.Range("L97").Value = 676.6667
$ws.Range("M97").Value = -87.33330000000001
$ws.Range("N97").Value = -1668.6667
$ws.Range("H116").Value = 2126.3333
$ws.Range("I116").Value = 2242.125
$ws.Range("J116").Value = 1200
$ws.Range("K116").Value = 2242.125
$ws.Range("L116").Value = 1200
$ws.Range("M116").Value = 51.875
$ws.Range("N116").Value = -5788
$ws.Range("H132").Value = 2823.75
$ws.Range("I132").Value = 1797
$ws.Range("J132").Value = 10011
$ws.Range("K132").Value = 5391
$ws.Range("L132").Value = 30033
$ws.Range("M132").Value = -2861
$ws.Range("N132").Value = -35093

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2126.3333
$ws.Range("I3").Value = 2242.125
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 2242.125
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = -2128.125
$ws.Range("N3").Value = -1428
$ws.Range("H64").Value = 1621.6666
$ws.Range("I64").Value = 1466.3334
$ws.Range("J64").Value = 1699.3334
$ws.Range("K64").Value = 1466.3334
$ws.Range("L64").Value = 1699.3334
$ws.Range("M64").Value = -1241.3334
$ws.Range("N64").Value = -2149.3334
$ws.Range("H67").Value = 1621.6666
$ws.Range("I67").Value = 1466.3334
$ws.Range("J67").Value = 1699.3334
$ws.Range("K67").Value = 1466.3334
$ws.Range("L67").Value = 1699.3334
$ws.Range("M67").Value = -686.3334
$ws.Range("N67").Value = -3259.3334
$ws.Range("H105").Value = 3017.1458
$ws.Range("I105").Value = 2447.1516
$ws.Range("J105").Value = 4271.1333
$ws.Range("K105").Value = 2447.1516
$ws.Range("L105").Value = 4271.1333
$ws.Range("M105").Value = -700.1516000000001
$ws.Range("N105").Value = -7765.1333
$ws.Range("H107").Value = 4300.864
$ws.Range("I107").Value = 2942.1875
$ws.Range("J107").Value = 7924
$ws.Range("K107").Value = 2942.1875
$ws.Range("L107").Value = 7924
$ws.Range("M107").Value = -1022.1875
$ws.Range("N107").Value = -11764
$ws.Range("H134").Value = 2650.5
$ws.Range("I134").Value = 1274.8
$ws.Range("J134").Value = 6089.75
$ws.Range("K134").Value = 3824.4
$ws.Range("L134").Value = 18269.25
$ws.Range("M134").Value = -1289.4
$ws.Range("N134").Value = -23339.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4911.8647
$ws.Range("I31").Value = 4539.4116
$ws.Range("J31").Value = 5228.45
$ws.Range("K31").Value = 4539.4116
$ws.Range("L31").Value = 5228.45
$ws.Range("M31").Value = -4244.4116
$ws.Range("N31").Value = -5818.45
$ws.Range("H34").Value = 4911.8647
$ws.Range("I34").Value = 4539.4116
$ws.Range("J34").Value = 5228.45
$ws.Range("K34").Value = 4539.4116
$ws.Range("L34").Value = 5228.45
$ws.Range("M34").Value = -4337.4116
$ws.Range("N34").Value = -5632.45
$ws.Range("H132").Value = 2100.7188
$ws.Range("I132").Value = 1955.7241
$ws.Range("J132").Value = 3502.3333
$ws.Range("K132").Value = 5867.1723
$ws.Range("L132").Value = 10506.9999
$ws.Range("M132").Value = -3337.1723
$ws.Range("N132").Value = -15566.9999
$ws.Range("H134").Value = 3495.647
$ws.Range("I134").Value = 3180.1
$ws.Range("J134").Value = 3946.4285
$ws.Range("K134").Value = 9540.299999999999
$ws.Range("L134").Value = 11839.2855
$ws.Range("M134").Value = -7005.299999999999
$ws.Range("N134").Value = -16909.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 8307.875
$ws.Range("I59").Value = 7554.3335
$ws.Range("J59").Value = 8760
$ws.Range("K59").Value = 22663.0005
$ws.Range("L59").Value = 26280
$ws.Range("M59").Value = -22123.0005
$ws.Range("N59").Value = -27360
$ws.Range("H60").Value = 793
$ws.Range("I60").Value = 858.25
$ws.Range("J60").Value = 749.5
$ws.Range("K60").Value = 2574.75
$ws.Range("L60").Value = 2248.5
$ws.Range("M60").Value = -2323.75
$ws.Range("N60").Value = -2750.5
$ws.Range("H81").Value = 3643
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 3643
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10929
$ws.Range("N81").Value = -13175
$ws.Range("H84").Value = 3643
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3643
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 32787
$ws.Range("N84").Value = -44019
$ws.Range("H140").Value = 4238.5884
$ws.Range("I140").Value = 3470.4666
$ws.Range("J140").Value = 9999.5
$ws.Range("K140").Value = 10411.3998
$ws.Range("L140").Value = 29998.5
$ws.Range("M140").Value = -5231.399800000001
$ws.Range("N140").Value = -40358.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6794.6
$ws.Range("I70").Value = 4988
$ws.Range("J70").Value = 7568.857
$ws.Range("K70").Value = 4988
$ws.Range("L70").Value = 7568.857
$ws.Range("M70").Value = -4718
$ws.Range("N70").Value = -8108.857
$ws.Range("H73").Value = 6794.6
$ws.Range("I73").Value = 4988
$ws.Range("J73").Value = 7568.857
$ws.Range("K73").Value = 4988
$ws.Range("L73").Value = 7568.857
$ws.Range("M73").Value = -4052
$ws.Range("N73").Value = -9440.857
$ws.Range("H107").Value = 92.75
$ws.Range("I107").Value = 95.75
$ws.Range("J107").Value = 89.75
$ws.Range("K107").Value = 95.75
$ws.Range("L107").Value = 89.75
$ws.Range("M107").Value = 1824.25
$ws.Range("N107").Value = -3929.75
$ws.Range("H113").Value = 4045.8572
$ws.Range("I113").Value = 1273.6666
$ws.Range("J113").Value = 4801.909
$ws.Range("K113").Value = 1273.6666
$ws.Range("L113").Value = 4801.909
$ws.Range("M113").Value = 896.3334
$ws.Range("N113").Value = -9141.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4850
$ws.Range("I7").Value = 4850
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4850
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -4738
$ws.Range("H40").Value = 3949.5
$ws.Range("I40").Value = 3949.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3949.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3813.5
$ws.Range("H82").Value = 3273
$ws.Range("I82").Value = 3437.4
$ws.Range("J82").Value = 2999
$ws.Range("K82").Value = 3437.4
$ws.Range("L82").Value = 2999
$ws.Range("M82").Value = -3076.4
$ws.Range("N82").Value = -3721
$ws.Range("H85").Value = 3273
$ws.Range("I85").Value = 3437.4
$ws.Range("J85").Value = 2999
$ws.Range("K85").Value = 3437.4
$ws.Range("L85").Value = 2999
$ws.Range("M85").Value = -2189.4
$ws.Range("N85").Value = -5495
$ws.Range("H100").Value = 2630.2
$ws.Range("I100").Value = 2787.75
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2787.75
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -2246.75
$ws.Range("N100").Value = -3082
$ws.Range("H126").Value = 4850
$ws.Range("I126").Value = 4850
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14550
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -12080
$ws.Range("H132").Value = 4923.091
$ws.Range("I132").Value = 3619
$ws.Range("J132").Value = 6227.1816
$ws.Range("K132").Value = 10857
$ws.Range("L132").Value = 18681.5448
$ws.Range("M132").Value = -8327
$ws.Range("N132").Value = -23741.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4250.75
$ws.Range("I81").Value = 4250.5
$ws.Range("J81").Value = 4251
$ws.Range("K81").Value = 8501
$ws.Range("L81").Value = 8502
$ws.Range("M81").Value = -7440
$ws.Range("N81").Value = -10624
$ws.Range("H84").Value = 4250.75
$ws.Range("I84").Value = 4250.5
$ws.Range("J84").Value = 4251
$ws.Range("K84").Value = 42505
$ws.Range("L84").Value = 42510
$ws.Range("M84").Value = -37201
$ws.Range("N84").Value = -53118
$ws.Range("H93").Value = 65000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 65000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 65000
$ws.Range("N93").Value = -69992
$ws.Range("H132").Value = 2820.5
$ws.Range("I132").Value = 2353.3
$ws.Range("J132").Value = 3599.1667
$ws.Range("K132").Value = 7059.900000000001
$ws.Range("L132").Value = 10797.5001
$ws.Range("M132").Value = -4529.900000000001
$ws.Range("N132").Value = -15857.5001
$ws.Range("H136").Value = 93748.63
$ws.Range("I136").Value = 1655.75
$ws.Range("J136").Value = 339329.66
$ws.Range("K136").Value = 4967.25
$ws.Range("L136").Value = 1017988.98
$ws.Range("M136").Value = -2417.25
$ws.Range("N136").Value = -1023088.98
